$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.533.25'
$ws.Range("D3").Value = '''3.672.52'
$ws.Range("E3").Value = '  -7.00%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''579.42'
$ws.Range("E5").Value = '  -4.04%  '
$ws.Range("D6").Value = '''169.59'
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D7").Value = '''3.665.31'
$ws.Range("E7").Value = '  -7.03%  '
$ws.Range("D8").Value = '''0.620'
$ws.Range("E8").Value = '  -9.07%  '
$ws.Range("D9").Value = '''0.998'
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("E10").Value = '  -11.21%  '
$ws.Range("E11").Value = '  -11.88%  '
$ws.Range("D12").Value = '''50.96'
$ws.Range("E12").Value = '  -9.18%  '
$ws.Range("D13").Value = '''0.0000285'
$ws.Range("E13").Value = '  -12.77%  '
$ws.Range("D14").Value = '''10.36'
$ws.Range("E14").Value = '  -10.91%  '
$ws.Range("D15").Value = '''4.243.32'
$ws.Range("E15").Value = '  -7.22%  '
$ws.Range("D16").Value = '''3.667.54'
$ws.Range("E16").Value = '  -7.50%  '
$ws.Range("D17").Value = '''19.26'
$ws.Range("E17").Value = '  -10.43%  '
$ws.Range("E18").Value = '  -3.55%  '
$ws.Range("E19").Value = '  -9.49%  '
$ws.Range("E20").Value = '  -9.92%  '
$ws.Range("D21").Value = '''67.342.08'
$ws.Range("E21").Value = '  -7.30%  '
$ws.Range("D22").Value = '''403.52'
$ws.Range("E22").Value = '  -9.36%  '
$ws.Range("D23").Value = '''4.48'
$ws.Range("E23").Value = '  -7.03%  '
$ws.Range("D24").Value = '''87.19'
$ws.Range("D25").Value = '''3.02'
$ws.Range("E25").Value = '  -9.44%  '
$ws.Range("E26").Value = '  -11.06%  '
$ws.Range("D27").Value = '''10.70'
$ws.Range("E27").Value = '  -5.22%  '
$ws.Range("D29").Value = '''3.76'
$ws.Range("E29").Value = '  -11.29%  '
$ws.Range("E30").Value = '  -10.05%  '
$ws.Range("D31").Value = '''32.33'
$ws.Range("E31").Value = '  -9.91%  '
$ws.Range("D32").Value = '''7.44'
$ws.Range("E32").Value = '  -5.86%  '
$ws.Range("D33").Value = '''12.32'
$ws.Range("E33").Value = '  -11.46%  '
$ws.Range("E34").Value = '  -10.19%  '
$ws.Range("D35").Value = '''64.33'
$ws.Range("E35").Value = '  -6.88%  '
$ws.Range("D36").Value = '''43.00'
$ws.Range("E36").Value = '  -13.36%  '
$ws.Range("D37").Value = '''589.90'
$ws.Range("E37").Value = '  -6.71%  '
$ws.Range("D38").Value = '''0.0₃0884'
$ws.Range("E38").Value = '  -10.93%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").Value = '''0.392'
$ws.Range("E40").Value = '  -8.61%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("E42").Value = '  -8.68%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = '''2.73'
$ws.Range("E43").Value = '  +3.61%  '
$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").Value = '''2.96'
$ws.Range("E44").Value = '  -13.53%  '
$ws.Range("E45").Value = '  -9.94%  '
$ws.Range("D46").Value = '''2.80'
$ws.Range("E46").Value = '  -13.09%  '
$ws.Range("D47").Value = '''9.09'
$ws.Range("E47").Value = '  -14.53%  '
$ws.Range("D48").Value = '''2.752.02'
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("D49").Value = '''0.132'
$ws.Range("E49").Value = '  -10.24%  '
$ws.Range("D50").Value = '''3.13'
$ws.Range("E50").Value = '  -7.38%  '
$ws.Range("D51").Value = '''2.66'
$ws.Range("E51").Value = '  -4.65%  '
